# Insert a new weekly price record as row 20, pushing the existing rows
# 20-32 down to 21-33 (dimension grows from A1:R32 to A1:R33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20..32 down by one row.
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new record.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44603
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = 100112028
$ws.Range("G20").Value = "Sandia"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Tercera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 280
$ws.Range("L20").Value = 300
$ws.Range("M20").Value = 290
$ws.Range("N20").Value = "$/kilo (volumen en unidades)"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 290
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"
